# Simulated Wild Card round and logged it
# Update "R" (road/playoff) row totals on both the OFF and DEF sheets
# to reflect the simulated Wild Card round game.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 ("R") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 284
$wsOff.Range("C3").Value = 201
$wsOff.Range("D3").Value = 55
$wsOff.Range("E3").Value = 25
$wsOff.Range("G3").Value = 6

# --- DEF sheet: row 3 ("R") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 192
$wsDef.Range("C3").Value = 135
$wsDef.Range("D3").Value = 47
$wsDef.Range("E3").Value = 24

$wb.Save()
